$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) keeps its text formatting so values like "1.00" or "6.50"
# are not auto-converted to numbers and lose formatting/trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '58.836.24'
$ws.Range('E2').Value = '  +3.59%  '
$ws.Range('D3').Value = '2.584.14'
$ws.Range('E3').Value = '  +2.06%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '519.67'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').Value = '140.06'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').Value = '  +1.33%  '
$ws.Range('D9').Value = '2.599.05'
$ws.Range('E9').Value = '  +2.27%  '
$ws.Range('D10').Value = '6.50'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('E11').Value = '  +1.96%  '
$ws.Range('D12').Value = '0.331'
$ws.Range('E12').Value = '  +3.05%  '
$ws.Range('E13').Value = '  +2.60%  '
$ws.Range('D14').Value = '3.052.32'
$ws.Range('D15').Value = '58.833.50'
$ws.Range('E15').Value = '  +3.55%  '
$ws.Range('D16').Value = '20.52'
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('D17').Value = '2.611.32'
$ws.Range('E17').Value = '  +2.85%  '
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').Value = '339.10'
$ws.Range('E19').Value = '  +2.85%  '
$ws.Range('D20').Value = '4.33'
$ws.Range('E20').Value = '  +2.11%  '
$ws.Range('D21').Value = '10.19'
$ws.Range('E21').Value = '  +1.65%  '
$ws.Range('D22').Value = '6.52'
$ws.Range('E22').Value = '  +6.60%  '
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').Value = '66.06'
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').Value = '0.405'
$ws.Range('E26').Value = '  +1.73%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = '7.11'
$ws.Range('E28').Value = '  +3.50%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').Value = '0.0₃0725'
$ws.Range('E30').Value = '  -1.96%  '
$ws.Range('D31').Value = '5.94'
$ws.Range('E31').Value = '  -4.44%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '18.77'
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.57'
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('D34').Value = '148.84'
$ws.Range('E34').Value = '  +0.39%  '
$ws.Range('E35').Value = '  +0.82%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').Value = '36.26'
$ws.Range('E37').Value = '  +1.76%  '
$ws.Range('D38').Value = '0.836'
$ws.Range('E38').Value = '  +2.35%  '
$ws.Range('E39').Value = '  +2.79%  '
$ws.Range('D40').Value = '0.824'
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('D41').Value = '3.51'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('D42').Value = '0.997'
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').Value = '275.73'
$ws.Range('E43').Value = '  +5.08%  '
$ws.Range('D44').Value = '10.74'
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('D45').Value = '0.0952'
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('D46').Value = '0.589'
$ws.Range('D47').Value = '0.0522'
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('D48').Value = '18.65'
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('D49').Value = '1.985.50'
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0220'
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '4.48'
$ws.Range('E51').Value = '  +0.03%  '
